# "Feat : calculate Bobot KHS"
# Renumber the NIM column, add new student rows (11-21), and apply the
# left/wrap text formatting used for the NIM column throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix up the existing rows (2-10): renumber NIM values and restyle A.
# ---------------------------------------------------------------------

# Row 2 & Row 5 keep the original (Arial) cell font, just gain left align.
$ws.Range("A2").Value = 9999999999
$ws.Range("A2").HorizontalAlignment = -4131

$ws.Range("A5").Value = 9999999993
$ws.Range("A5").HorizontalAlignment = -4131

# Row 3 & 4 switch back to the default (Calibri) font, with left align.
$ws.Range("A3").ClearFormats()
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").Value = 9999999991

$ws.Range("A4").ClearFormats()
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").Value = 9999999992

# Rows 6-10 switch to default font, left align, vertical-center + wrap.
$ws.Range("A6").ClearFormats()
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true

$ws.Range("A7").ClearFormats()
$ws.Range("A7").HorizontalAlignment = -4131
$ws.Range("A7").VerticalAlignment = -4108
$ws.Range("A7").WrapText = $true

$ws.Range("A8").ClearFormats()
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("A8").WrapText = $true

$ws.Range("A9").ClearFormats()
$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A9").WrapText = $true

$ws.Range("A10").ClearFormats()
$ws.Range("A10").HorizontalAlignment = -4131
$ws.Range("A10").VerticalAlignment = -4108
$ws.Range("A10").WrapText = $true

# ---------------------------------------------------------------------
# 2) Append the new rows (11-21) with the continuing NIM sequence.
# ---------------------------------------------------------------------

# Row 11 : NIM stored as TEXT (9999999910) with the "@" number format.
$ws.Range("A11").ClearFormats()
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").HorizontalAlignment = -4131
$ws.Range("A11").VerticalAlignment = -4108
$ws.Range("A11").Value = "9999999910"
$ws.Range("B11").Value = "satri"
$ws.Range("C11").Value = 90
$ws.Range("D11").Value = 100
$ws.Range("E11").Value = 88
$ws.Range("F11").Value = 100
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 50
$ws.Range("I11").Value = 88

$ws.Range("A12").ClearFormats()
$ws.Range("A12").HorizontalAlignment = -4131
$ws.Range("A12").VerticalAlignment = -4108
$ws.Range("A12").WrapText = $true
$ws.Range("A12").Value = 9999999911
$ws.Range("B12").Value = "fajr"
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = 88
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 90
$ws.Range("G12").Value = 88
$ws.Range("H12").Value = 88
$ws.Range("I12").Value = 88

$ws.Range("A13").ClearFormats()
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("A13").WrapText = $true
$ws.Range("A13").Value = 9999999912
$ws.Range("B13").Value = "absa"
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = 90
$ws.Range("E13").Value = 90
$ws.Range("F13").Value = 100
$ws.Range("G13").Value = 100
$ws.Range("H13").Value = 90
$ws.Range("I13").Value = 100

$ws.Range("A14").ClearFormats()
$ws.Range("A14").HorizontalAlignment = -4131
$ws.Range("A14").VerticalAlignment = -4108
$ws.Range("A14").WrapText = $true
$ws.Range("A14").Value = 9999999913
$ws.Range("B14").Value = "han"
$ws.Range("C14").Value = 90
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 88
$ws.Range("G14").Value = 88
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 90

$ws.Range("A15").ClearFormats()
$ws.Range("A15").HorizontalAlignment = -4131
$ws.Range("A15").VerticalAlignment = -4108
$ws.Range("A15").WrapText = $true
$ws.Range("A15").Value = 9999999914
$ws.Range("B15").Value = "das"
$ws.Range("C15").Value = 90
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = 88
$ws.Range("F15").Value = 100
$ws.Range("G15").Value = 50
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 88

$ws.Range("A16").ClearFormats()
$ws.Range("A16").HorizontalAlignment = -4131
$ws.Range("A16").VerticalAlignment = -4108
$ws.Range("A16").WrapText = $true
$ws.Range("A16").Value = 9999999915
$ws.Range("B16").Value = "eq"
$ws.Range("C16").Value = 50
$ws.Range("D16").Value = 88
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 90
$ws.Range("G16").Value = 88
$ws.Range("H16").Value = 88
$ws.Range("I16").Value = 88

$ws.Range("A17").ClearFormats()
$ws.Range("A17").HorizontalAlignment = -4131
$ws.Range("A17").VerticalAlignment = -4108
$ws.Range("A17").WrapText = $true
$ws.Range("A17").Value = 9999999916
$ws.Range("B17").Value = "fsda"
$ws.Range("C17").Value = 50
$ws.Range("D17").Value = 90
$ws.Range("E17").Value = 90
$ws.Range("F17").Value = 100
$ws.Range("G17").Value = 100
$ws.Range("H17").Value = 90
$ws.Range("I17").Value = 100

$ws.Range("A18").ClearFormats()
$ws.Range("A18").HorizontalAlignment = -4131
$ws.Range("A18").VerticalAlignment = -4108
$ws.Range("A18").WrapText = $true
$ws.Range("A18").Value = 99999999917
$ws.Range("B18").Value = "dsd"
$ws.Range("C18").Value = 90
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 88
$ws.Range("G18").Value = 88
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 90

$ws.Range("A19").ClearFormats()
$ws.Range("A19").HorizontalAlignment = -4131
$ws.Range("A19").VerticalAlignment = -4108
$ws.Range("A19").WrapText = $true
$ws.Range("A19").Value = 9999999918
$ws.Range("B19").Value = "sdewq"
$ws.Range("C19").Value = 90
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 88
$ws.Range("G19").Value = 88
$ws.Range("H19").Value = 100
$ws.Range("I19").Value = 90

$ws.Range("A20").ClearFormats()
$ws.Range("A20").HorizontalAlignment = -4131
$ws.Range("A20").VerticalAlignment = -4108
$ws.Range("A20").WrapText = $true
$ws.Range("A20").Value = 9999999919
$ws.Range("B20").Value = "satri"
$ws.Range("C20").Value = 90
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = 88
$ws.Range("F20").Value = 100
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 88

# Row 21 gets an explicit black font color (distinct from the theme color
# used elsewhere), left align, vertical-center and wrap.
$ws.Range("A21").ClearFormats()
$ws.Range("A21").Font.Color = 0
$ws.Range("A21").HorizontalAlignment = -4131
$ws.Range("A21").VerticalAlignment = -4108
$ws.Range("A21").WrapText = $true
$ws.Range("A21").Value = 9999999920
$ws.Range("B21").Value = "fajr"
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 88
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 88
$ws.Range("I21").Value = 88

# ---------------------------------------------------------------------
# 3) Selection / active cell, matching the recorded sheet view.
# ---------------------------------------------------------------------
$ws.Range("K12").Select()
